$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'maa://24702 (94.59), maa://25390 (96.25), maa://36681 (87.34)'
$ws.Range("L2").Value = '*maa://24633 (56.17), *maa://30515 (70.48), maa://39402 (92.75), *maa://34787 (72.73), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range("T2").Value = 'maa://22742 (90.7), *maa://20791 (62.34)'
$ws.Range("AB2").Value = 'maa://21246 (91.44), maa://36684 (95.83), ***maa://22731 (6.25)'
$ws.Range("AF2").Value = 'maa://25251 (91.74), ***maa://21730 (27.27), ***maa://39501 (15.15), **maa://36675 (50.0)'
$ws.Range("D3").Value = 'maa://40192 (96.67), maa://36987 (96.15), maa://39849 (88.89)'
$ws.Range("L3").Value = '*maa://22880 (64.85), maa://20276 (86.67), *maa://22749 (76.92)'
$ws.Range("P3").Value = 'maa://21249 (94.65), maa://26254 (96.67), **maa://22738 (50.0)'
$ws.Range("T3").Value = 'maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.18), maa://45854 (84.85)'
$ws.Range("X3").Value = 'maa://27396 (84.36), maa://27484 (96.67), maa://27480 (83.33)'
$ws.Range("D4").Value = 'maa://24632 (93.71), **maa://24303 (38.46), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range("T4").Value = 'maa://32509 (95.8), maa://27295 (86.49), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range("X4").Value = '**maa://32495 (48.89), ***maa://31785 (22.22), maa://43217 (92.0), ***maa://36683 (28.26)'
$ws.Range("D5").Value = 'maa://21245 (84.49), maa://22744 (84.62)'
$ws.Range("D6").Value = 'maa://42407 (96.72)'
$ws.Range("T6").Value = 'maa://37411 (88.24)'
$ws.Range("X7").Value = 'maa://22399 (95.62), *maa://22758 (75.0)'
$ws.Range("A8").Value = '更新日期：2025.03.28 13:20:26'
$ws.Range("D8").Value = '*maa://21476 (73.58), *maa://39431 (60.0), *maa://37551 (57.14)'
$ws.Range("X8").Value = 'maa://21411 (96.0)'
$ws.Range("T9").Value = '**maa://22866 (30.19), maa://26222 (98.15)'
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '6'
$ws.Range("AB9").Value = 'maa://28711 (87.2), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), *maa://45044 (66.67), maa://40166 (96.3)'
$ws.Range("AF9").Value = 'maa://26206 (88.19), *maa://22865 (51.85)'
$ws.Range("D10").Value = '***maa://25695 (18.32), ***maa://39951 (13.56), ***maa://34206 (22.22), ***maa://39243 (25.0), *maa://45271 (53.49)'
$ws.Range("X10").Value = 'maa://22301 (97.78), maa://45828 (87.5), maa://22726 (100.0)'
$ws.Range("AF10").Value = '*maa://25021 (54.0), *maa://22733 (62.16), **maa://22761 (50.0)'
$ws.Range("X11").Value = 'maa://36713 (97.81)'
$ws.Range("H12").Value = 'maa://21867 (90.06), ***maa://45826 (25.0)'
$ws.Range("AB12").Value = 'maa://23669 (95.52), maa://36677 (94.12), maa://39872 (92.0)'
$ws.Range("AF12").Value = '*maa://28932 (78.15), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range("D13").Value = 'maa://24999 (92.14), maa://36673 (92.41), maa://25001 (85.92)'
$ws.Range("H13").Value = '*maa://21248 (74.09), **maa://22728 (46.67)'
$ws.Range("X13").Value = 'maa://34957 (81.71), **maa://22768 (50.0)'
$ws.Range("AF13").Value = '**maa://22737 (34.25), maa://39883 (90.79), *maa://39885 (53.33)'
$ws.Range("L14").Value = 'maa://26245 (96.77), maa://21288 (96.3), maa://39841 (94.02), maa://36682 (97.44)'
$ws.Range("P14").Value = 'maa://23250 (98.76), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range("T14").Value = 'maa://22521 (94.44), maa://42751 (100.0)'
$ws.Range("H15").Value = 'maa://24304 (87.95), maa://21478 (89.47)'
$ws.Range("P15").Value = 'maa://24762 (90.53), *maa://22727 (70.0)'
$ws.Range("T15").Value = 'maa://23892 (96.3)'
$ws.Range("AF15").Value = 'maa://21364 (81.18), *maa://36666 (77.59), *maa://22766 (68.33)'
$ws.Range("D16").Value = 'maa://21441 (96.43), maa://36679 (94.55), maa://37650 (97.67)'
$ws.Range("T16").Value = 'maa://22729 (94.55), *maa://28648 (69.57), *maa://36674 (79.63)'
$ws.Range("H17").Value = 'maa://22430 (88.83), maa://39599 (83.64)'
$ws.Range("D18").Value = 'maa://24570 (96.98)'
$ws.Range("L18").Value = 'maa://22466 (90.75), *maa://22732 (51.04)'
$ws.Range("X18").Value = 'maa://21917 (96.97), maa://22741 (87.5)'
$ws.Range("T19").Value = 'maa://24386 (99.18)'
$ws.Range("D20").Value = 'maa://21432 (90.56), maa://25198 (93.69), *maa://20795 (50.77), maa://36680 (91.18)'
$ws.Range("H20").Value = 'maa://22864 (90.0)'
$ws.Range("P20").Value = 'maa://37442 (95.35)'
$ws.Range("AF21").Value = 'maa://22524 (93.39), *maa://22432 (78.31)'
$ws.Range("T22").Value = 'maa://38495 (87.5)'
$ws.Range("X22").Value = 'maa://21282 (98.63), *maa://37649 (65.52)'
$ws.Range("P23").Value = 'maa://30587 (92.0), *maa://29748 (76.15), ***maa://29785 (16.18), *maa://37566 (77.5)'
$ws.Range("D24").Value = '*maa://24368 (78.61), *maa://46650 (62.5)'
$ws.Range("X24").Value = 'maa://29988 (83.97), maa://23504 (93.33), **maa://22892 (40.54), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)'
$ws.Range("AF24").Value = 'maa://22523 (85.22), *maa://36672 (79.31), maa://29910 (93.22), **maa://21440 (35.71), maa://45831 (85.71)'
$ws.Range("H25").Value = '*maa://29063 (72.62), *maa://25311 (74.77), ***maa://22725 (4.76), *maa://45047 (66.67)'
$ws.Range("AB25").Value = 'maa://31215 (88.33), maa://24516 (80.22), maa://26001 (87.5)'
$ws.Range("D26").Value = 'maa://41802 (90.91)'
$ws.Range("AB26").Value = 'maa://42235 (94.69)'
$ws.Range("X28").Value = 'maa://39929 (90.89), maa://41749 (91.4), ***maa://39723 (13.89)'
$ws.Range("AF28").Value = 'maa://36660 (92.58), *maa://36701 (66.67)'
$ws.Range("H29").Value = '*maa://25175 (66.67)'
$ws.Range("L29").Value = 'maa://28432 (93.59), maa://28440 (80.87), maa://31400 (98.82), *maa://28650 (71.43)'
$ws.Range("P29").Value = '*maa://23168 (58.06), *maa://30050 (55.56)'
$ws.Range("AF29").Value = '*maa://24080 (68.93), maa://42865 (81.16), ***maa://34960 (8.33)'
$ws.Range("D30").Value = 'maa://45792 (94.12)'
$ws.Range("H32").Value = 'maa://21895 (97.07), maa://36667 (97.73), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range("T32").Value = 'maa://42859 (96.15), maa://41108 (88.0), maa://41238 (97.2), maa://45523 (100.0)'
$ws.Range("P34").Value = 'maa://48817 (90.91)'
$ws.Range("L35").Value = 'maa://41296 (96.47)'
$ws.Range("L37").Value = 'maa://45718 (97.6), *maa://47069 (73.33), maa://45789 (100.0)'
$ws.Range("H39").Value = 'maa://36670 (89.22), maa://25199 (84.82), maa://30434 (92.05), maa://45059 (81.82), ***maa://25036 (19.23), *maa://44165 (66.67)'
$ws.Range("P40").Value = 'maa://23278 (95.31), maa://21386 (95.79), maa://36664 (89.29), maa://45550 (87.5)'
$ws.Range("H41").Value = 'maa://24466 (93.88)'
$ws.Range("H44").Value = 'maa://29768 (98.08), maa://27728 (96.12)'
$ws.Range("H46").Value = 'maa://35931 (91.98), maa://43901 (93.33)'
$ws.Range("H47").Value = 'maa://27410 (96.5), maa://29661 (97.35), maa://28038 (84.62)'
$ws.Range("H53").Value = 'maa://32534 (94.22), **maa://32434 (33.33)'
$ws.Range("H55").Value = 'maa://32532 (92.11)'
$ws.Range("H57").Value = 'maa://25176 (98.44)'
$ws.Range("H58").Value = '*maa://37964 (58.14)'
$ws.Range("H59").Value = 'maa://31270 (94.78), maa://27746 (82.46)'
